# Auto-generated Excel COM-interop script applying the diff to Sargatanas_Profits workbook
# Each worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) corresponds to one table block in the diff.
$wb = $excel.ActiveWorkbook

# ---------- Sheet: ALC ----------
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J46").Value = 5000
$ws.Range("L46").Value = 15000
$ws.Range("N46").Value = -15238
$ws.Range("H55").Value = 421.66666
$ws.Range("I55").Value = 700
$ws.Range("K55").Value = 700
$ws.Range("M55").Value = -486
$ws.Range("J60").Value = 5000
$ws.Range("L60").Value = 15000
$ws.Range("N60").Value = -15968
$ws.Range("H70").Value = 53033420
$ws.Range("I70").Value = 62502676
$ws.Range("J70").Value = 47622420
$ws.Range("K70").Value = 187508028
$ws.Range("L70").Value = 142867260
$ws.Range("M70").Value = -187507758
$ws.Range("N70").Value = -142867800
$ws.Range("H73").Value = 53033420
$ws.Range("I73").Value = 62502676
$ws.Range("J73").Value = 47622420
$ws.Range("K73").Value = 187508028
$ws.Range("L73").Value = 142867260
$ws.Range("M73").Value = -187507092
$ws.Range("N73").Value = -142869132
$ws.Range("H80").Value = 34776.535
$ws.Range("I80").Value = 14920.143
$ws.Range("J80").Value = 52150.875
$ws.Range("K80").Value = 44760.429
$ws.Range("L80").Value = 156452.625
$ws.Range("M80").Value = -43762.429
$ws.Range("N80").Value = -158448.625
$ws.Range("H83").Value = 34776.535
$ws.Range("I83").Value = 14920.143
$ws.Range("J83").Value = 52150.875
$ws.Range("K83").Value = 134281.287
$ws.Range("L83").Value = 469357.875
$ws.Range("M83").Value = -129289.287
$ws.Range("N83").Value = -479341.875
$ws.Range("H86").Value = 62138940
$ws.Range("I86").Value = 93753340
$ws.Range("J86").Value = 7942821
$ws.Range("K86").Value = 93753340
$ws.Range("L86").Value = 7942821
$ws.Range("M86").Value = -93752217
$ws.Range("N86").Value = -7945067
$ws.Range("H88").Value = 22270844
$ws.Range("I88").Value = 66670156
$ws.Range("J88").Value = 71186
$ws.Range("K88").Value = 66670156
$ws.Range("L88").Value = 71186
$ws.Range("M88").Value = -66669750
$ws.Range("N88").Value = -71998
$ws.Range("H89").Value = 62138940
$ws.Range("I89").Value = 93753340
$ws.Range("J89").Value = 7942821
$ws.Range("K89").Value = 468766700
$ws.Range("L89").Value = 39714105
$ws.Range("M89").Value = -468761084
$ws.Range("N89").Value = -39725337
$ws.Range("H91").Value = 22270844
$ws.Range("I91").Value = 66670156
$ws.Range("J91").Value = 71186
$ws.Range("K91").Value = 66670156
$ws.Range("L91").Value = 71186
$ws.Range("M91").Value = -66668752
$ws.Range("N91").Value = -73994
$ws.Range("H101").Value = 599.7
$ws.Range("I101").Value = 387.85715
$ws.Range("J101").Value = 1094
$ws.Range("K101").Value = 1163.57145
$ws.Range("L101").Value = 3282
$ws.Range("M101").Value = 458.4285500000001
$ws.Range("N101").Value = -6526
$ws.Range("H103").Value = 1280.0667
$ws.Range("I103").Value = 614.5
$ws.Range("J103").Value = 1382.4615
$ws.Range("K103").Value = 1843.5
$ws.Range("L103").Value = 4147.3845
$ws.Range("M103").Value = -1257.5
$ws.Range("N103").Value = -5319.3845
$ws.Range("H111").Value = 125000000
$ws.Range("I111").Value = 125000000
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 375000000
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -374996933
$ws.Range("N111").ClearContents()
$ws.Range("H123").Value = 51000
$ws.Range("J123").Value = 51000
$ws.Range("L123").Value = 51000
$ws.Range("N123").Value = -60800
$ws.Range("H132").Value = 1450.9474
$ws.Range("I132").Value = 1403.5714
$ws.Range("J132").Value = 2003.6666
$ws.Range("K132").Value = 4210.7142
$ws.Range("L132").Value = 6010.9998
$ws.Range("M132").Value = -1680.7142
$ws.Range("N132").Value = -11070.9998
$ws.Range("H137").Value = 2193.8276
$ws.Range("I137").Value = 2286.8572
$ws.Range("J137").Value = 1949.625
$ws.Range("K137").Value = 6860.571599999999
$ws.Range("L137").Value = 5848.875
$ws.Range("M137").Value = -4310.571599999999
$ws.Range("N137").Value = -10948.875

# ---------- Sheet: ARM ----------
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1408806.9
$ws.Range("I32").Value = 1491790.5
$ws.Range("J32").Value = 14681.2
$ws.Range("K32").Value = 1491790.5
$ws.Range("L32").Value = 14681.2
$ws.Range("M32").Value = -1491503.5
$ws.Range("N32").Value = -15255.2
$ws.Range("I34").Value = 84500
$ws.Range("J34").Value = 499999.5
$ws.Range("K34").Value = 84500
$ws.Range("L34").Value = 499999.5
$ws.Range("M34").Value = -84229
$ws.Range("N34").Value = -500541.5
$ws.Range("H45").Value = 15990.667
$ws.Range("I45").Value = 6250.5
$ws.Range("K45").Value = 6250.5
$ws.Range("M45").Value = -5873.5
$ws.Range("H61").Value = 2309.4868
$ws.Range("I61").Value = 1508.5374
$ws.Range("K61").Value = 1508.5374
$ws.Range("M61").Value = -1296.5374
$ws.Range("H62").Value = 43000
$ws.Range("J62").Value = 43000
$ws.Range("L62").Value = 43000
$ws.Range("N62").Value = -44248
$ws.Range("H63").Value = 2496.5
$ws.Range("I63").Value = 2496.5
$ws.Range("K63").Value = 2496.5
$ws.Range("M63").Value = -1810.5
$ws.Range("H65").Value = 43000
$ws.Range("J65").Value = 43000
$ws.Range("L65").Value = 129000
$ws.Range("N65").Value = -135240
$ws.Range("H66").Value = 2496.5
$ws.Range("I66").Value = 2496.5
$ws.Range("K66").Value = 12482.5
$ws.Range("M66").Value = -9050.5
$ws.Range("H74").Value = 45945.55
$ws.Range("I74").Value = 69530.03999999999
$ws.Range("J74").Value = 5515
$ws.Range("K74").Value = 69530.03999999999
$ws.Range("L74").Value = 5515
$ws.Range("M74").Value = -68656.03999999999
$ws.Range("N74").Value = -7263
$ws.Range("H77").Value = 45945.55
$ws.Range("I77").Value = 69530.03999999999
$ws.Range("J77").Value = 5515
$ws.Range("K77").Value = 347650.2
$ws.Range("L77").Value = 27575
$ws.Range("M77").Value = -343282.2
$ws.Range("N77").Value = -36311
$ws.Range("H88").Value = 1829
$ws.Range("J88").Value = 1928.8572
$ws.Range("L88").Value = 1928.8572
$ws.Range("N88").Value = -2740.8572
$ws.Range("H91").Value = 1829
$ws.Range("J91").Value = 1928.8572
$ws.Range("L91").Value = 1928.8572
$ws.Range("N91").Value = -4736.8572
$ws.Range("H97").Value = 16688963
$ws.Range("I97").Value = 434
$ws.Range("J97").Value = 27814648
$ws.Range("K97").Value = 434
$ws.Range("L97").Value = 27814648
$ws.Range("M97").Value = 62
$ws.Range("N97").Value = -27815640
$ws.Range("H110").Value = 23810654
$ws.Range("I110").Value = 1073.5454
$ws.Range("K110").Value = 1073.5454
$ws.Range("M110").Value = 971.4546
$ws.Range("H119").Value = 67971.664
$ws.Range("J119").Value = 67971.664
$ws.Range("L119").Value = 67971.664
$ws.Range("N119").Value = -77647.664
$ws.Range("H122").Value = 22408.545
$ws.Range("I122").Value = 26811.75
$ws.Range("K122").Value = 80435.25
$ws.Range("M122").Value = -77985.25
$ws.Range("H132").Value = 3282.111
$ws.Range("I132").Value = 2019.3077
$ws.Range("J132").Value = 6565.4
$ws.Range("K132").Value = 6057.9231
$ws.Range("L132").Value = 19696.2
$ws.Range("M132").Value = -3527.9231
$ws.Range("N132").Value = -24756.2
$ws.Range("H136").Value = 2309.4868
$ws.Range("I136").Value = 1508.5374
$ws.Range("K136").Value = 4525.6122
$ws.Range("M136").Value = -1975.6122

# ---------- Sheet: BSM ----------
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 12822751
$ws.Range("I20").Value = 33336106
$ws.Range("J20").Value = 1904.125
$ws.Range("K20").Value = 33336106
$ws.Range("L20").Value = 1904.125
$ws.Range("M20").Value = -33335859
$ws.Range("N20").Value = -2398.125
$ws.Range("H86").Value = 31285252
$ws.Range("I86").Value = 50999.57
$ws.Range("J86").Value = 90914280
$ws.Range("K86").Value = 50999.57
$ws.Range("L86").Value = 90914280
$ws.Range("M86").Value = -49876.57
$ws.Range("N86").Value = -90916526
$ws.Range("H89").Value = 31285252
$ws.Range("I89").Value = 50999.57
$ws.Range("J89").Value = 90914280
$ws.Range("K89").Value = 254997.85
$ws.Range("L89").Value = 454571400
$ws.Range("M89").Value = -249381.85
$ws.Range("N89").Value = -454582632
$ws.Range("H94").Value = 1443.4762
$ws.Range("I94").Value = 555.86664
$ws.Range("K94").Value = 555.86664
$ws.Range("M94").Value = -104.86664
$ws.Range("H99").Value = 9095226
$ws.Range("I99").Value = 4739.7144
$ws.Range("J99").Value = 30306360
$ws.Range("K99").Value = 4739.7144
$ws.Range("L99").Value = 30306360
$ws.Range("M99").Value = -3241.7144
$ws.Range("N99").Value = -30309356
$ws.Range("H134").Value = 3296.9424
$ws.Range("I134").Value = 1121.3334
$ws.Range("J134").Value = 7075.6313
$ws.Range("K134").Value = 3364.0002
$ws.Range("L134").Value = 21226.8939
$ws.Range("M134").Value = -829.0001999999999
$ws.Range("N134").Value = -26296.8939

# ---------- Sheet: CRP ----------
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3781.3333
$ws.Range("I16").Value = 2142
$ws.Range("K16").Value = 2142
$ws.Range("M16").Value = -1855
$ws.Range("H22").Value = 290.5
$ws.Range("I22").Value = 304.2857
$ws.Range("J22").Value = 258.33334
$ws.Range("K22").Value = 304.2857
$ws.Range("L22").Value = 258.33334
$ws.Range("M22").Value = 45.71429999999998
$ws.Range("N22").Value = -958.33334
$ws.Range("H31").Value = 6656.347
$ws.Range("I31").Value = 2316.9285
$ws.Range("J31").Value = 12442.238
$ws.Range("K31").Value = 2316.9285
$ws.Range("L31").Value = 12442.238
$ws.Range("M31").Value = -2021.9285
$ws.Range("N31").Value = -13032.238
$ws.Range("H34").Value = 6656.347
$ws.Range("I34").Value = 2316.9285
$ws.Range("J34").Value = 12442.238
$ws.Range("K34").Value = 2316.9285
$ws.Range("L34").Value = 12442.238
$ws.Range("M34").Value = -2114.9285
$ws.Range("N34").Value = -12846.238
$ws.Range("H58").Value = 10421915
$ws.Range("I58").Value = 17858660
$ws.Range("J58").Value = 10471.15
$ws.Range("K58").Value = 17858660
$ws.Range("L58").Value = 10471.15
$ws.Range("M58").Value = -17858457
$ws.Range("N58").Value = -10877.15
$ws.Range("H59").Value = 67500
$ws.Range("I59").Value = 35000
$ws.Range("K59").Value = 35000
$ws.Range("M59").Value = -33855
$ws.Range("H62").Value = 13894996
$ws.Range("I62").Value = 25006292
$ws.Range("J62").Value = 5876
$ws.Range("K62").Value = 25006292
$ws.Range("L62").Value = 5876
$ws.Range("M62").Value = -25005668
$ws.Range("N62").Value = -7124
$ws.Range("H65").Value = 13894996
$ws.Range("I65").Value = 25006292
$ws.Range("J65").Value = 5876
$ws.Range("K65").Value = 125031460
$ws.Range("L65").Value = 29380
$ws.Range("M65").Value = -125028340
$ws.Range("N65").Value = -35620
$ws.Range("H76").Value = 4953.75
$ws.Range("I76").Value = 4953.75
$ws.Range("K76").Value = 4953.75
$ws.Range("M76").Value = -4638.75
$ws.Range("H79").Value = 4953.75
$ws.Range("I79").Value = 4953.75
$ws.Range("K79").Value = 4953.75
$ws.Range("M79").Value = -3861.75
$ws.Range("H86").Value = 7880711
$ws.Range("I86").Value = 15631922
$ws.Range("K86").Value = 15631922
$ws.Range("M86").Value = -15630799
$ws.Range("H89").Value = 7880711
$ws.Range("I89").Value = 15631922
$ws.Range("K89").Value = 78159610
$ws.Range("M89").Value = -78153994
$ws.Range("H99").Value = 10192.462
$ws.Range("I99").Value = 11417.333
$ws.Range("J99").Value = 9142.571
$ws.Range("K99").Value = 11417.333
$ws.Range("L99").Value = 9142.571
$ws.Range("M99").Value = -9919.333000000001
$ws.Range("N99").Value = -12138.571
$ws.Range("H107").Value = 3351.85
$ws.Range("I107").Value = 3440.9
$ws.Range("J107").Value = 3262.8
$ws.Range("K107").Value = 3440.9
$ws.Range("L107").Value = 3262.8
$ws.Range("M107").Value = -1520.9
$ws.Range("N107").Value = -7102.8
$ws.Range("H113").Value = 3781.3333
$ws.Range("I113").Value = 2142
$ws.Range("K113").Value = 2142
$ws.Range("M113").Value = 28
$ws.Range("H126").Value = 10192.462
$ws.Range("I126").Value = 11417.333
$ws.Range("J126").Value = 9142.571
$ws.Range("K126").Value = 34251.999
$ws.Range("L126").Value = 27427.713
$ws.Range("M126").Value = -31781.999
$ws.Range("N126").Value = -32367.713
$ws.Range("H132").Value = 4281.9214
$ws.Range("I132").Value = 1748.3334
$ws.Range("J132").Value = 10362.533
$ws.Range("K132").Value = 5245.0002
$ws.Range("L132").Value = 31087.599
$ws.Range("M132").Value = -2715.0002
$ws.Range("N132").Value = -36147.599
$ws.Range("H134").Value = 6284.486
$ws.Range("I134").Value = 1436.5
$ws.Range("J134").Value = 8813.869000000001
$ws.Range("K134").Value = 4309.5
$ws.Range("L134").Value = 26441.607
$ws.Range("M134").Value = -1774.5
$ws.Range("N134").Value = -31511.607
$ws.Range("H136").Value = 10421915
$ws.Range("I136").Value = 17858660
$ws.Range("J136").Value = 10471.15
$ws.Range("K136").Value = 53575980
$ws.Range("L136").Value = 31413.45
$ws.Range("M136").Value = -53573430
$ws.Range("N136").Value = -36513.45
$ws.Range("H141").Value = 59698.047
$ws.Range("J141").Value = 60631.477
$ws.Range("L141").Value = 60631.477
$ws.Range("N141").Value = -70991.477

# ---------- Sheet: CUL ----------
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1278608.8
$ws.Range("I11").Value = 1438348.6
$ws.Range("J11").Value = 690
$ws.Range("K11").Value = 4315045.800000001
$ws.Range("L11").Value = 2070
$ws.Range("M11").Value = -4314905.800000001
$ws.Range("N11").Value = -2350
$ws.Range("H38").Value = 70.5
$ws.Range("J38").Value = 78
$ws.Range("L38").Value = 234
$ws.Range("N38").Value = -928
$ws.Range("H56").Value = 7242
$ws.Range("I56").Value = 7242
$ws.Range("K56").Value = 7242
$ws.Range("M56").Value = -6712
$ws.Range("H87").Value = 1000000000
$ws.Range("I87").Value = 1000000000
$ws.Range("K87").Value = 3000000000
$ws.Range("M87").Value = -2999998752
$ws.Range("H90").Value = 1000000000
$ws.Range("I90").Value = 1000000000
$ws.Range("K90").Value = 9000000000
$ws.Range("M90").Value = -8999993760
$ws.Range("H93").Value = 6477.778
$ws.Range("J93").Value = 6975
$ws.Range("L93").Value = 20925
$ws.Range("N93").Value = -24669
$ws.Range("H122").Value = 4716538
$ws.Range("I122").Value = 9429843
$ws.Range("K122").Value = 84868587
$ws.Range("M122").Value = -84866137
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H131").Value = 998.4167
$ws.Range("I131").Value = 747.1053000000001
$ws.Range("J131").Value = 1953.4
$ws.Range("K131").Value = 2241.3159
$ws.Range("L131").Value = 5860.200000000001
$ws.Range("M131").Value = 2798.6841
$ws.Range("N131").Value = -15940.2

# ---------- Sheet: GSM ----------
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 49.066666
$ws.Range("I2").Value = 41.75
$ws.Range("K2").Value = 41.75
$ws.Range("M2").Value = 71.25
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("H48").Value = 20000
$ws.Range("J48").Value = 20000
$ws.Range("L48").Value = 20000
$ws.Range("N48").Value = -20970
$ws.Range("H70").Value = 8945.556
$ws.Range("I70").Value = 7400.7144
$ws.Range("J70").Value = 9928.637000000001
$ws.Range("K70").Value = 7400.7144
$ws.Range("L70").Value = 9928.637000000001
$ws.Range("M70").Value = -7130.7144
$ws.Range("N70").Value = -10468.637
$ws.Range("H73").Value = 8945.556
$ws.Range("I73").Value = 7400.7144
$ws.Range("J73").Value = 9928.637000000001
$ws.Range("K73").Value = 7400.7144
$ws.Range("L73").Value = 9928.637000000001
$ws.Range("M73").Value = -6464.7144
$ws.Range("N73").Value = -11800.637
$ws.Range("H80").Value = 4748.3335
$ws.Range("J80").Value = 4623.25
$ws.Range("L80").Value = 4623.25
$ws.Range("N80").Value = -6619.25
$ws.Range("H83").Value = 4748.3335
$ws.Range("J83").Value = 4623.25
$ws.Range("L83").Value = 23116.25
$ws.Range("N83").Value = -33100.25
$ws.Range("H96").Value = 53101.5
$ws.Range("J96").Value = 53101.5
$ws.Range("L96").Value = 53101.5
$ws.Range("N96").Value = -58593.5
$ws.Range("H102").Value = 3180.5264
$ws.Range("I102").Value = 3314
$ws.Range("K102").Value = 3314
$ws.Range("M102").Value = -1692
$ws.Range("H113").Value = 6778.484
$ws.Range("I113").Value = 3250
$ws.Range("K113").Value = 3250
$ws.Range("M113").Value = -1080
$ws.Range("H122").Value = 4530826
$ws.Range("I122").Value = 6039102
$ws.Range("J122").Value = 5998.5
$ws.Range("K122").Value = 18117306
$ws.Range("L122").Value = 17995.5
$ws.Range("M122").Value = -18114856
$ws.Range("N122").Value = -22895.5
$ws.Range("H126").Value = 2643.4546
$ws.Range("I126").Value = 2428.875
$ws.Range("J126").Value = 3215.6667
$ws.Range("K126").Value = 7286.625
$ws.Range("L126").Value = 9647.000100000001
$ws.Range("M126").Value = -4816.625
$ws.Range("N126").Value = -14587.0001
$ws.Range("H132").Value = 1651.5714
$ws.Range("I132").Value = 1399.1333
$ws.Range("J132").Value = 3166.2
$ws.Range("K132").Value = 4197.3999
$ws.Range("L132").Value = 9498.599999999999
$ws.Range("M132").Value = -1667.3999
$ws.Range("N132").Value = -14558.6

# ---------- Sheet: LTW ----------
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5532.25
$ws.Range("I7").Value = 4244.9
$ws.Range("J7").Value = 8750.625
$ws.Range("K7").Value = 4244.9
$ws.Range("L7").Value = 8750.625
$ws.Range("M7").Value = -4132.9
$ws.Range("N7").Value = -8974.625
$ws.Range("H22").Value = 1639.2632
$ws.Range("J22").Value = 3199.3333
$ws.Range("L22").Value = 3199.3333
$ws.Range("N22").Value = -3789.3333
$ws.Range("H27").Value = 1639.2632
$ws.Range("J27").Value = 3199.3333
$ws.Range("L27").Value = 3199.3333
$ws.Range("N27").Value = -3413.3333
$ws.Range("H40").Value = 6862.5625
$ws.Range("I40").Value = 5310.8887
$ws.Range("K40").Value = 5310.8887
$ws.Range("M40").Value = -5174.8887
$ws.Range("H46").Value = 1918074.1
$ws.Range("I46").Value = 11494770
$ws.Range("J46").Value = 2734.8667
$ws.Range("K46").Value = 11494770
$ws.Range("L46").Value = 2734.8667
$ws.Range("M46").Value = -11494582
$ws.Range("N46").Value = -3110.8667
$ws.Range("H55").Value = 30303340
$ws.Range("I55").Value = 62500084
$ws.Range("J55").Value = 521.2353000000001
$ws.Range("K55").Value = 62500084
$ws.Range("L55").Value = 521.2353000000001
$ws.Range("M55").Value = -62499911
$ws.Range("N55").Value = -867.2353000000001
$ws.Range("H61").Value = 9032.777
$ws.Range("I61").Value = 7800
$ws.Range("J61").Value = 9186.875
$ws.Range("K61").Value = 7800
$ws.Range("L61").Value = 9186.875
$ws.Range("M61").Value = -7598
$ws.Range("N61").Value = -9590.875
$ws.Range("H68").Value = 2817
$ws.Range("I68").Value = 1952.9286
$ws.Range("K68").Value = 1952.9286
$ws.Range("M68").Value = -1203.9286
$ws.Range("H71").Value = 2817
$ws.Range("I71").Value = 1952.9286
$ws.Range("K71").Value = 9764.643
$ws.Range("M71").Value = -6020.643
$ws.Range("H93").Value = 5618.067
$ws.Range("I93").Value = 3296.7
$ws.Range("J93").Value = 10260.8
$ws.Range("K93").Value = 3296.7
$ws.Range("L93").Value = 10260.8
$ws.Range("M93").Value = -2048.7
$ws.Range("N93").Value = -12756.8
$ws.Range("H106").Value = 25714.6
$ws.Range("J106").Value = 25714.6
$ws.Range("L106").Value = 25714.6
$ws.Range("N106").Value = -28238.6
$ws.Range("H113").Value = 9032.777
$ws.Range("I113").Value = 7800
$ws.Range("J113").Value = 9186.875
$ws.Range("K113").Value = 7800
$ws.Range("L113").Value = 9186.875
$ws.Range("M113").Value = -5630
$ws.Range("N113").Value = -13526.875
$ws.Range("H122").Value = 3920.82
$ws.Range("I122").Value = 2866.1082
$ws.Range("K122").Value = 8598.3246
$ws.Range("M122").Value = -6148.3246
$ws.Range("H126").Value = 5532.25
$ws.Range("I126").Value = 4244.9
$ws.Range("J126").Value = 8750.625
$ws.Range("K126").Value = 12734.7
$ws.Range("L126").Value = 26251.875
$ws.Range("M126").Value = -10264.7
$ws.Range("N126").Value = -31191.875
$ws.Range("H132").Value = 10874475
$ws.Range("I132").Value = 17244340
$ws.Range("J132").Value = 8235.235000000001
$ws.Range("K132").Value = 51733020
$ws.Range("L132").Value = 24705.705
$ws.Range("M132").Value = -51730490
$ws.Range("N132").Value = -29765.705
$ws.Range("H136").Value = 6111.24
$ws.Range("I136").Value = 2309
$ws.Range("J136").Value = 8250
$ws.Range("K136").Value = 6927
$ws.Range("L136").Value = 24750
$ws.Range("M136").Value = -4377
$ws.Range("N136").Value = -29850

# ---------- Sheet: WVR ----------
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H81").Value = 11117003
$ws.Range("I81").Value = 2948.4285
$ws.Range("J81").Value = 50016196
$ws.Range("K81").Value = 5896.857
$ws.Range("L81").Value = 100032392
$ws.Range("M81").Value = -4835.857
$ws.Range("N81").Value = -100034514
$ws.Range("H84").Value = 11117003
$ws.Range("I84").Value = 2948.4285
$ws.Range("J84").Value = 50016196
$ws.Range("K84").Value = 29484.285
$ws.Range("L84").Value = 500161960
$ws.Range("M84").Value = -24180.285
$ws.Range("N84").Value = -500172568
$ws.Range("H107").Value = 12346336
$ws.Range("J107").Value = 30303910
$ws.Range("L107").Value = 90911730
$ws.Range("N107").Value = -90915570
$ws.Range("H113").Value = 1528.5
$ws.Range("I113").Value = 1270.7142
$ws.Range("K113").Value = 3812.1426
$ws.Range("M113").Value = -1642.1426
$ws.Range("H119").Value = 56319
$ws.Range("J119").Value = 56319
$ws.Range("L119").Value = 56319
$ws.Range("N119").Value = -65995
$ws.Range("H120").Value = 54189
$ws.Range("J120").Value = 54189
$ws.Range("L120").Value = 54189
$ws.Range("N120").Value = -63865
$ws.Range("H121").Value = 54459.332
$ws.Range("J121").Value = 54459.332
$ws.Range("L121").Value = 54459.332
$ws.Range("N121").Value = -57953.332
$ws.Range("H122").Value = 105710.64
$ws.Range("I122").Value = 135132.67
$ws.Range("J122").Value = 7637.222
$ws.Range("K122").Value = 405398.01
$ws.Range("L122").Value = 22911.666
$ws.Range("M122").Value = -402948.01
$ws.Range("N122").Value = -27811.666
$ws.Range("H126").Value = 2555.2222
$ws.Range("I126").Value = 999.7143
$ws.Range("K126").Value = 2999.1429
$ws.Range("M126").Value = -529.1428999999998
$ws.Range("H132").Value = 11637695
$ws.Range("I132").Value = 16671019
$ws.Range("J132").Value = 22331.46
$ws.Range("K132").Value = 50013057
$ws.Range("L132").Value = 66994.38
$ws.Range("M132").Value = -50010527
$ws.Range("N132").Value = -72054.38
$ws.Range("H136").Value = 18541010
$ws.Range("I136").Value = 29412890
$ws.Range("J136").Value = 58815.55
$ws.Range("K136").Value = 88238670
$ws.Range("L136").Value = 176446.65
$ws.Range("M136").Value = -88236120
$ws.Range("N136").Value = -181546.65
$ws.Range("H137").Value = 58638.75
$ws.Range("J137").Value = 58638.75
$ws.Range("L137").Value = 58638.75
$ws.Range("N137").Value = -68838.75
$ws.Range("H138").Value = 80000
$ws.Range("J138").Value = 80000
$ws.Range("L138").Value = 80000
$ws.Range("N138").Value = -90280
$ws.Range("H139").Value = 88813.71000000001
$ws.Range("J139").Value = 88616
$ws.Range("L139").Value = 88616
$ws.Range("N139").Value = -98896
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 100000
$ws.Range("J141").Value = 100000
$ws.Range("L141").Value = 100000
$ws.Range("N141").Value = -110360
